# Fix loi ko ton tai maKhuyenMai (trang ban sach)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the book titles (remove stray trailing index numbers / fix ordering)
$ws.Range("B3").Value = "Conan"
$ws.Range("B4").Value = "Năm mươi Sắc thái"
$ws.Range("B5").Value = "Cho tôi một vé đi tuổi thơ"
$ws.Range("B6").Value = "Chiến Thắng Con Quỷ Trong Bạn"
$ws.Range("B7").Value = "Đất Rừng Phương Nam"
$ws.Range("B8").Value = "Harry Potter Và Hòn Đá Phù Thuỷ"

# Update quantities to 50 for each book row
$ws.Range("C3").Value = 50
$ws.Range("C4").Value = 50
$ws.Range("C5").Value = 50
$ws.Range("C6").Value = 50
$ws.Range("C7").Value = 50
$ws.Range("C8").Value = 50

# Update the active selection to match the saved view state
$ws.Range("H7").Select()

$wb.Save()
